# Auto-generated edit script: updates Leve profit-calculation cells (H-N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, per scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 41.75
$ws.Range("I39").Value = 19.2
$ws.Range("K39").Value = 57.59999999999999
$ws.Range("M39").Value = 238.4
$ws.Range("H51").Value = 2720.4119
$ws.Range("H137").Value = 1298.2285
$ws.Range("I137").Value = 1078.92
$ws.Range("J137").Value = 1846.5
$ws.Range("K137").Value = 3236.76
$ws.Range("L137").Value = 5539.5
$ws.Range("M137").Value = -686.7600000000002
$ws.Range("N137").Value = -10639.5
$ws.Range("H141").Value = 836.2857
$ws.Range("I141").Value = 836.2857
$ws.Range("K141").Value = 2508.8571
$ws.Range("M141").Value = 2671.1429

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9300.816000000001
$ws.Range("I32").Value = 5529.122
$ws.Range("J32").Value = 28630.75
$ws.Range("K32").Value = 5529.122
$ws.Range("L32").Value = 28630.75
$ws.Range("M32").Value = -5242.122
$ws.Range("N32").Value = -29204.75
$ws.Range("H43").Value = 44495
$ws.Range("J43").Value = 44495
$ws.Range("L43").Value = 44495
$ws.Range("N43").Value = -45121
$ws.Range("H61").Value = 3544.7705
$ws.Range("I61").Value = 2452.925
$ws.Range("K61").Value = 2452.925
$ws.Range("M61").Value = -2240.925
$ws.Range("H113").Value = 30398
$ws.Range("J113").Value = 30398
$ws.Range("L113").Value = 30398
$ws.Range("N113").Value = -39076
$ws.Range("H125").Value = 26928.334
$ws.Range("J125").Value = 26928.334
$ws.Range("L125").Value = 26928.334
$ws.Range("N125").Value = -36768.334
$ws.Range("H132").Value = 3187.608
$ws.Range("I132").Value = 2786.2046
$ws.Range("J132").Value = 5710.7144
$ws.Range("K132").Value = 8358.613799999999
$ws.Range("L132").Value = 17132.1432
$ws.Range("M132").Value = -5828.613799999999
$ws.Range("N132").Value = -22192.1432
$ws.Range("H135").Value = 39142.855
$ws.Range("J135").Value = 39142.855
$ws.Range("L135").Value = 39142.855
$ws.Range("N135").Value = -49282.855
$ws.Range("H136").Value = 3544.7705
$ws.Range("I136").Value = 2452.925
$ws.Range("K136").Value = 7358.775000000001
$ws.Range("M136").Value = -4808.775000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 18494
$ws.Range("I82").Value = 18494
$ws.Range("K82").Value = 18494
$ws.Range("M82").Value = -18111
$ws.Range("H85").Value = 18494
$ws.Range("I85").Value = 18494
$ws.Range("K85").Value = 18494
$ws.Range("M85").Value = -17168
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 575.5
$ws.Range("I10").Value = 167.33333
$ws.Range("J10").Value = 1800
$ws.Range("K10").Value = 167.33333
$ws.Range("L10").Value = 1800
$ws.Range("M10").Value = -28.33332999999999
$ws.Range("N10").Value = -2078
$ws.Range("H19").Value = 385.5
$ws.Range("I19").Value = 385.5
$ws.Range("K19").Value = 385.5
$ws.Range("M19").Value = -215.5
$ws.Range("H24").Value = 385.5
$ws.Range("I24").Value = 385.5
$ws.Range("K24").Value = 385.5
$ws.Range("M24").Value = -215.5
$ws.Range("H55").Value = 44000
$ws.Range("J55").Value = 44000
$ws.Range("L55").Value = 44000
$ws.Range("N55").Value = -44630
$ws.Range("H122").Value = 2590.1667
$ws.Range("J122").Value = 3435
$ws.Range("L122").Value = 10305
$ws.Range("N122").Value = -15205
$ws.Range("H134").Value = 12909.8125
$ws.Range("I134").Value = 5101.8716
$ws.Range("K134").Value = 15305.6148
$ws.Range("M134").Value = -12770.6148

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 1707
$ws.Range("I8").Value = 1707
$ws.Range("K8").Value = 5121
$ws.Range("M8").Value = -4982
$ws.Range("H87").Value = 22224.875
$ws.Range("J87").Value = 29997.5
$ws.Range("L87").Value = 89992.5
$ws.Range("N87").Value = -92488.5
$ws.Range("H90").Value = 22224.875
$ws.Range("J90").Value = 29997.5
$ws.Range("L90").Value = 269977.5
$ws.Range("N90").Value = -282457.5
$ws.Range("H99").Value = 1279.4
$ws.Range("I99").Value = 1279.4
$ws.Range("K99").Value = 3838.2
$ws.Range("M99").Value = -1592.2
$ws.Range("H107").Value = 829.2
$ws.Range("I107").Value = 661.5
$ws.Range("J107").Value = 1500
$ws.Range("K107").Value = 1984.5
$ws.Range("L107").Value = 4500
$ws.Range("M107").Value = -64.5
$ws.Range("N107").Value = -8340
$ws.Range("H122").Value = 2996.4285
$ws.Range("J122").Value = 5000
$ws.Range("L122").Value = 45000
$ws.Range("N122").Value = -49900

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 15753.765
$ws.Range("I126").Value = 20743.416
$ws.Range("J126").Value = 3778.6
$ws.Range("K126").Value = 62230.24800000001
$ws.Range("L126").Value = 11335.8
$ws.Range("M126").Value = -59760.24800000001
$ws.Range("N126").Value = -16275.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2835.9333
$ws.Range("I22").Value = 2744
$ws.Range("J22").Value = 3019.8
$ws.Range("K22").Value = 2744
$ws.Range("L22").Value = 3019.8
$ws.Range("M22").Value = -2449
$ws.Range("N22").Value = -3609.8
$ws.Range("H27").Value = 2835.9333
$ws.Range("I27").Value = 2744
$ws.Range("J27").Value = 3019.8
$ws.Range("K27").Value = 2744
$ws.Range("L27").Value = 3019.8
$ws.Range("M27").Value = -2637
$ws.Range("N27").Value = -3233.8
$ws.Range("H40").Value = 3730.625
$ws.Range("I40").Value = 2723.8333
$ws.Range("J40").Value = 6751
$ws.Range("K40").Value = 2723.8333
$ws.Range("L40").Value = 6751
$ws.Range("M40").Value = -2587.8333
$ws.Range("N40").Value = -7023
$ws.Range("H122").Value = 283528.38
$ws.Range("I122").Value = 375513.12
$ws.Range("J122").Value = 7574.1113
$ws.Range("K122").Value = 1126539.36
$ws.Range("L122").Value = 22722.3339
$ws.Range("M122").Value = -1124089.36
$ws.Range("N122").Value = -27622.3339
$ws.Range("H140").Value = 231175.6
$ws.Range("J140").Value = 271499.75
$ws.Range("L140").Value = 271499.75
$ws.Range("N140").Value = -281859.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 12998.5
$ws.Range("J15").Value = 18997
$ws.Range("L15").Value = 18997
$ws.Range("N15").Value = -19573
$ws.Range("H58").Value = 45999.5
$ws.Range("J58").Value = 45999.5
$ws.Range("L58").Value = 45999.5
$ws.Range("N58").Value = -46615.5
$ws.Range("H126").Value = 12633.182
$ws.Range("I126").Value = 19329.334
$ws.Range("J126").Value = 4597.8
$ws.Range("K126").Value = 57988.00199999999
$ws.Range("L126").Value = 13793.4
$ws.Range("M126").Value = -55518.00199999999
$ws.Range("N126").Value = -18733.4
$ws.Range("H131").Value = 144999.5
$ws.Range("J131").Value = 144999.5
$ws.Range("L131").Value = 144999.5
$ws.Range("N131").Value = -155079.5
$ws.Range("H132").Value = 1536.6666
$ws.Range("I132").Value = 1486.2059
$ws.Range("K132").Value = 4458.6177
$ws.Range("M132").Value = -1928.6177
